$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.TrimEnd([char]13) -eq $needle) {
            return $p
        }
    }
    return $null
}

# =====================================================================
# 1. "Dokumentacja" bullet: small word-level edits (tworzyliśmy -> tworzono,
#    umieściliśmy -> umieszczono)
# =====================================================================
$d.Content.Find.Execute("tworzyliśmy przez cały czas pracy nad aplikacją", $true, $false, $false, $false, $false, `
    $true, 1, $false, "tworzono przez cały czas pracy nad aplikacją", 2) | Out-Null

$d.Content.Find.Execute("umieściliśmy komentarze", $true, $false, $false, $false, $false, `
    $true, 1, $false, "umieszczono komentarze", 2) | Out-Null

# =====================================================================
# 2. Two new paragraphs after the "Dokumentacja" bullet:
#       a) a new bulleted ("Akapitzlist"/numId 3) paragraph about
#          differences vs. the initial plan
#       b) a plain justified paragraph in red: "tutaj analiza SWAT"
# =====================================================================
$nbsp = [char]0xA0
$pDokNeedle = "Dokumentacja – dokumentację tworzono przez cały czas pracy nad aplikacją. Pierwszym jej etapem było uzasadnienie biznesowe, następnie określiliśmy wymagania aplikacji. Dodatkowo umieszczono komentarze w kodzie, a także instrukcję dla użytkowników. Ostatnim etapem tworzenia dokumentacji jest niniejszy raport. Nad dokumentacją pracowali wszyscy członkowie zespołu: uzasadnienie biznesowe stworzyliśmy wspólnie, podobnie jak wymagania. Natomiast komentarze w kodzie i instrukcję dla użytkownika napisali Kinga Dobrowolska i Błażej Kurzep, zaś raport Monika Czajka i Zuzanna Kontna. Średnio każdy z" + $nbsp + "członków poświęcił na to 5 godzin."
$pDok = Find-ParagraphByText $d $pDokNeedle
if ($pDok -eq $null) { throw "Could not locate the 'Dokumentacja' paragraph after the word edits." }

$pDok.Range.InsertParagraphAfter()
$pRoznice = $pDok.Next()
$pRoznice.Range.Text = "Różnice w porównaniu z planem początkowym – pierwotny plan nie zakładał pojawienia się instrukcji obsługi wywoływanej z poziomu okna aplikacji, jednak zespołowi udało się to zrobić w ostatecznej wersji programu."

$pRoznice.Range.InsertParagraphAfter()
$pSwat = $pRoznice.Next()
$pSwat.Range.ListFormat.RemoveNumbers()
$pSwat.Style = "Normal"
$pSwat.Alignment = 3
$pSwat.Range.Text = "tutaj analiza SWAT"
$pSwat.Range.Font.Color = 255

# =====================================================================
# 3. Replace the first empty paragraph following "Opis uzyskanego
#    rezultatu" with three new content paragraphs.
# =====================================================================
$pResult = Find-ParagraphByText $d "Opis uzyskanego rezultatu"
if ($pResult -eq $null) { throw "Could not locate the 'Opis uzyskanego rezultatu' heading paragraph." }
$pEmpty = $pResult.Next()             # first empty <w:p/> to be filled in

$pEmpty.Range.Text = "Jak zostało wspomniane powyżej, największą zmianą w porównaniu do początkowej specyfikacji wymagań było pojawienie się instrukcji obsługi aplikacji dla użytkowników. Instrukcja dostępna jest z poziomu paska zadań programu i otwiera się w nowym oknie w postaci dokumentu z rozszerzeniem PDF. "
$pEmpty.FirstLineIndent = 18
$pEmpty.Alignment = 3

$pEmpty.Range.InsertParagraphAfter()
$pUseA = $pEmpty.Next()
$pUseA.Range.Text = "`tWszystkie wstępne założenia programu zostały w pełni zrealizowane. Zespół sprostał zarówno wymaganiom funkcjonalnym jaki i niefunkcjonalnym aplikacji. Aplikacja poprawnie oblicza wysokość zysku z lokaty oraz należny podatek przy zadanej kwocie początkowej, oprocentowaniu oraz okresie kapitalizacji odsetek. Jednym z głównych celów zespołu było stworzenie programu prostego w obsłudze, co ma swoje odzwierciedlenie w prostym i schludnym oknie aplikacji, a także dostępności instrukcji obsługi. "

$pUseA.Range.InsertParagraphAfter()
$pUseB = $pUseA.Next()
$pUseB.Range.Text = "`tKolejnym etapem ulepszenia oprogramowania mogłaby być opcja zmiany języka aplikacji, w szczególności na język angielski. Z czasem możliwe byłoby dodanie kolejnych opcji językowych. Kolejnym udogodnieniem mogła by stać się opcja możliwości wybrania odsetek przed końcem lokaty. Obie te aktualizacje wymagałyby ogromnego nakładu pracy ze strony zespołu, jednak mogłyby przyczynić się do większej funkcjonalności aplikacji, a dodatkowo zwiększyć grono odbiorców."
